$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Add a new sub-bullet paragraph after "Feature Engineering"
[void]$tr.InsertAfter("`rSustainabile vs Unsustainable")

$newPara = $tr.Paragraphs($tr.Paragraphs().Count)
$newPara.IndentLevel = 2
$newPara.Font.Size = 22

# Force the new paragraph's text to be split into two runs, matching
# how PowerPoint splits runs at a flagged (misspelled) word boundary.
$firstRun = $newPara.Characters(1, 12)
$firstRun.Text = "Sustainabile"
